# Update repo data based on source changes
#
# 1. Bump the dataset version number (package_description!D7): 0.3.0 -> 0.3.1
# 2. Refresh the per-local-authority lookup columns (region / type / county-la)
#    on the "promises" sheet. A handful of local authorities were abolished
#    in the 2023 local government reorganisation (Cumbria, North Yorkshire
#    and Somerset district councils merged into new unitary authorities), so
#    those codes no longer resolve against the refreshed reference data and
#    their region/type/county-la lookup values are cleared.
# 3. The "promises_metadata" sheet's F column no longer lists the full set of
#    valid options for the region/type columns, so those two cells are
#    cleared and the column is narrowed back down.

$wb = $excel.ActiveWorkbook

# --- 1. Bump version number on package_description -------------------------
$wsDesc = $wb.Worksheets.Item("package_description")
$wsDesc.Range("D7").Value2 = "0.3.1"

# --- 2. Clear stale region/type/county-la lookups on promises --------------
$wsPromises = $wb.Worksheets.Item("promises")

# Local-authority codes that were abolished/merged away and no longer have a
# match in the refreshed region/type/county-la reference data.
$staleCodes = @("ALL","BAR","CAR","COP","CRA","CMA","EDN","HAE","HAG","MEN","NYK","RIH","RYE","SCE","SEG","SEL","SOM","SLA","SSO","SWT")

$lastRow = $wsPromises.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $wsPromises.Cells.Item($r, 1).Value2
    if ($staleCodes -contains $code) {
        $wsPromises.Cells.Item($r, 10).Value2 = ""   # J: region
        $wsPromises.Cells.Item($r, 11).Value2 = ""   # K: local-authority-type-name
        $wsPromises.Cells.Item($r, 12).Value2 = ""   # L: county-la
    }
}

# --- 3. Clear the options list from promises_metadata and narrow column F --
$wsMeta = $wb.Worksheets.Item("promises_metadata")
$wsMeta.Range("F11").Value2 = ""
$wsMeta.Range("F12").Value2 = ""
$wsMeta.Columns.Item(6).ColumnWidth = 61.8
